$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3 & 4: "Lankadinee" -> "Rahul" person's data (name + email),
# the "state" (B) column values are left untouched (2 and 3 respectively).
$ws.Range("A3").Value = "Rahul"
$ws.Range("C3").Value = "rahulkalubowila@gmail.com"
$ws.Range("A4").Value = "Rahul"
$ws.Range("C4").Value = "rahulkalubowila@gmail.com"

# Drop the old row 5 (the separate "Rahul" entry that is now folded into
# rows 3 & 4 above).
$ws.Rows("5").Delete()

# Rebuild the hyperlinks: the sheet-wide hyperlink collection needs to go
# from { C2, C3:C4, C5(stale) } to { C2, C3, C4 }. Clear everything and
# re-add the three links so each data row carries its own single-cell link.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:lankadinee@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:lankadinee@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:rahulkalubowila@gmail.com")

# Restore the Hyperlink cell style (re-adding the links above reset it to a
# freshly minted duplicate style) so the cells keep using the original
# Hyperlink style.
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"

# Match the final selection left behind in the saved file.
$ws.Range("B4").Select()
